# Auto-generated update: refresh market price / profit figures across sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 18: You Grow, Girl / Growth Formula Beta
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).ClearContents()
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Cells.Item(28, 8).Value = 568.9
$ws.Cells.Item(28, 9).Value = 528.0952
$ws.Cells.Item(28, 10).Value = 664.1111
$ws.Cells.Item(28, 11).Value = 528.0952
$ws.Cells.Item(28, 12).Value = 664.1111
$ws.Cells.Item(28, 13).Value = -43.09519999999998
$ws.Cells.Item(28, 14).Value = -1634.1111
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Cells.Item(98, 8).Value = 1208.6111
$ws.Cells.Item(98, 9).Value = 1187.9166
$ws.Cells.Item(98, 10).Value = 1250
$ws.Cells.Item(98, 11).Value = 1187.9166
$ws.Cells.Item(98, 12).Value = 1250
$ws.Cells.Item(98, 13).Value = 310.0834
$ws.Cells.Item(98, 14).Value = -4246
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Cells.Item(122, 8).Value = 1208.6111
$ws.Cells.Item(122, 9).Value = 1187.9166
$ws.Cells.Item(122, 10).Value = 1250
$ws.Cells.Item(122, 11).Value = 3563.7498
$ws.Cells.Item(122, 12).Value = 3750
$ws.Cells.Item(122, 13).Value = -1113.7498
$ws.Cells.Item(122, 14).Value = -8650
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Cells.Item(135, 8).Value = 894.1053000000001
$ws.Cells.Item(135, 9).Value = 811.75
$ws.Cells.Item(135, 10).Value = 1333.3334
$ws.Cells.Item(135, 11).Value = 7305.75
$ws.Cells.Item(135, 12).Value = 12000.0006
$ws.Cells.Item(135, 13).Value = -4770.75
$ws.Cells.Item(135, 14).Value = -17070.0006
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 2669.3845
$ws.Cells.Item(137, 9).Value = 1857
$ws.Cells.Item(137, 10).Value = 3617.1667
$ws.Cells.Item(137, 11).Value = 5571
$ws.Cells.Item(137, 12).Value = 10851.5001
$ws.Cells.Item(137, 13).Value = -3021
$ws.Cells.Item(137, 14).Value = -15951.5001
# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Cells.Item(141, 8).Value = 11910.647
$ws.Cells.Item(141, 9).Value = 2425.8572
$ws.Cells.Item(141, 10).Value = 18550
$ws.Cells.Item(141, 11).Value = 7277.571599999999
$ws.Cells.Item(141, 12).Value = 55650
$ws.Cells.Item(141, 13).Value = -2097.571599999999
$ws.Cells.Item(141, 14).Value = -66010

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Cells.Item(102, 8).Value = 1940.3889
$ws.Cells.Item(102, 9).Value = 2013.8
$ws.Cells.Item(102, 10).Value = 1848.625
$ws.Cells.Item(102, 11).Value = 2013.8
$ws.Cells.Item(102, 12).Value = 1848.625
$ws.Cells.Item(102, 13).Value = -391.8
$ws.Cells.Item(102, 14).Value = -5092.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Cells.Item(20, 8).Value = 14291.6
$ws.Cells.Item(20, 9).Value = 4971.143
$ws.Cells.Item(20, 10).Value = 36039.332
$ws.Cells.Item(20, 11).Value = 4971.143
$ws.Cells.Item(20, 12).Value = 36039.332
$ws.Cells.Item(20, 13).Value = -4724.143
$ws.Cells.Item(20, 14).Value = -36533.332
# Row 80: Unbreaker / Titanium Ingot
$ws.Cells.Item(80, 8).Value = 3260.9524
$ws.Cells.Item(80, 9).Value = 799
$ws.Cells.Item(80, 10).Value = 5969.1
$ws.Cells.Item(80, 11).Value = 799
$ws.Cells.Item(80, 12).Value = 5969.1
$ws.Cells.Item(80, 13).Value = 199
$ws.Cells.Item(80, 14).Value = -7965.1
# Row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Cells.Item(83, 8).Value = 3260.9524
$ws.Cells.Item(83, 9).Value = 799
$ws.Cells.Item(83, 10).Value = 5969.1
$ws.Cells.Item(83, 11).Value = 3995
$ws.Cells.Item(83, 12).Value = 29845.5
$ws.Cells.Item(83, 13).Value = 997
$ws.Cells.Item(83, 14).Value = -39829.5
# Row 94: High Steal / High Steel Nugget
$ws.Cells.Item(94, 8).Value = 770.2
$ws.Cells.Item(94, 9).Value = 714.0476
$ws.Cells.Item(94, 10).Value = 1065
$ws.Cells.Item(94, 11).Value = 714.0476
$ws.Cells.Item(94, 12).Value = 1065
$ws.Cells.Item(94, 13).Value = -263.0476
$ws.Cells.Item(94, 14).Value = -1967
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 1847.909
$ws.Cells.Item(134, 9).Value = 796.25
$ws.Cells.Item(134, 10).Value = 4652.3335
$ws.Cells.Item(134, 11).Value = 2388.75
$ws.Cells.Item(134, 12).Value = 13957.0005
$ws.Cells.Item(134, 13).Value = 146.25
$ws.Cells.Item(134, 14).Value = -19027.0005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 50: The Arsenal of Theocracy / Cobalt Halberd
$ws.Cells.Item(50, 8).Value = 36922.223
$ws.Cells.Item(50, 9).Value = 8500
$ws.Cells.Item(50, 10).Value = 40475
$ws.Cells.Item(50, 11).Value = 8500
$ws.Cells.Item(50, 12).Value = 40475
$ws.Cells.Item(50, 13).Value = -7875
$ws.Cells.Item(50, 14).Value = -41725
# Row 51: Greenstone for Greenhorns / Jade Crook
$ws.Cells.Item(51, 8).Value = 25481.883
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 25481.883
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 25481.883
$ws.Cells.Item(51, 14).Value = -26953.883
# Row 59: Bow Down to Magic / Crab Bow
$ws.Cells.Item(59, 8).Value = 37644.445
$ws.Cells.Item(59, 9).Value = 4000
$ws.Cells.Item(59, 10).Value = 41850
$ws.Cells.Item(59, 11).Value = 4000
$ws.Cells.Item(59, 12).Value = 41850
$ws.Cells.Item(59, 13).Value = -2855
$ws.Cells.Item(59, 14).Value = -44140
# Row 60: Bowing to Greater Power / Yew Longbow
$ws.Cells.Item(60, 8).Value = 15457.695
$ws.Cells.Item(60, 9).Value = 15550
$ws.Cells.Item(60, 10).Value = 15448.904
$ws.Cells.Item(60, 11).Value = 15550
$ws.Cells.Item(60, 12).Value = 15448.904
$ws.Cells.Item(60, 13).Value = -15039
$ws.Cells.Item(60, 14).Value = -16470.904
# Row 61: Incant Now, Think Later / Jade Crook
$ws.Cells.Item(61, 8).Value = 25481.883
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 10).Value = 25481.883
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 12).Value = 25481.883
$ws.Cells.Item(61, 14).Value = -26177.883
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 5203.4546
$ws.Cells.Item(132, 9).Value = 5748.6665
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 17245.9995
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -14715.9995
$ws.Cells.Item(132, 14).Value = -20057

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 26: A Grape Idea / Grape Juice
$ws.Cells.Item(26, 8).Value = 683.2308
$ws.Cells.Item(26, 9).Value = 312
$ws.Cells.Item(26, 10).Value = 1001.4286
$ws.Cells.Item(26, 11).Value = 936
$ws.Cells.Item(26, 12).Value = 3004.2858
$ws.Cells.Item(26, 13).Value = -648
$ws.Cells.Item(26, 14).Value = -3580.2858
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 2365.6875
$ws.Cells.Item(131, 9).Value = 553.5
$ws.Cells.Item(131, 10).Value = 2412.1538
$ws.Cells.Item(131, 11).Value = 1660.5
$ws.Cells.Item(131, 12).Value = 7236.4614
$ws.Cells.Item(131, 13).Value = 3379.5
$ws.Cells.Item(131, 14).Value = -17316.4614

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Cells.Item(80, 8).Value = 2557.5
$ws.Cells.Item(80, 9).Value = 2708
$ws.Cells.Item(80, 10).Value = 2407
$ws.Cells.Item(80, 11).Value = 2708
$ws.Cells.Item(80, 12).Value = 2407
$ws.Cells.Item(80, 13).Value = -1710
$ws.Cells.Item(80, 14).Value = -4403
# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Cells.Item(83, 8).Value = 2557.5
$ws.Cells.Item(83, 9).Value = 2708
$ws.Cells.Item(83, 10).Value = 2407
$ws.Cells.Item(83, 11).Value = 13540
$ws.Cells.Item(83, 12).Value = 12035
$ws.Cells.Item(83, 13).Value = -8548
$ws.Cells.Item(83, 14).Value = -22019
# Row 94: Wants and Needles / Bombfish Needle
$ws.Cells.Item(94, 8).Value = 50000
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 50000
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 50000
$ws.Cells.Item(94, 14).Value = -51352
# Row 106: Choker in the Clutch / Palladium Choker of Aiming
$ws.Cells.Item(106, 8).Value = 29499.5
$ws.Cells.Item(106, 9).Value = 29000
$ws.Cells.Item(106, 10).Value = 29999
$ws.Cells.Item(106, 11).Value = 29000
$ws.Cells.Item(106, 12).Value = 29999
$ws.Cells.Item(106, 13).Value = -27738
$ws.Cells.Item(106, 14).Value = -32523
# Row 125: Pewter-hewn Punishment / Pewter Choker of Slaying
$ws.Cells.Item(125, 8).Value = 40000
$ws.Cells.Item(125, 9).Value = 0
$ws.Cells.Item(125, 10).Value = 40000
$ws.Cells.Item(125, 11).Value = 0
$ws.Cells.Item(125, 12).Value = 40000
$ws.Cells.Item(125, 14).Value = -44920
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 1852.5333
$ws.Cells.Item(126, 9).Value = 1566.9354
$ws.Cells.Item(126, 10).Value = 2484.9285
$ws.Cells.Item(126, 11).Value = 4700.8062
$ws.Cells.Item(126, 12).Value = 7454.7855
$ws.Cells.Item(126, 13).Value = -2230.8062
$ws.Cells.Item(126, 14).Value = -12394.7855

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Cells.Item(22, 8).Value = 2256.889
$ws.Cells.Item(22, 9).Value = 305
$ws.Cells.Item(22, 10).Value = 2814.5715
$ws.Cells.Item(22, 11).Value = 305
$ws.Cells.Item(22, 12).Value = 2814.5715
$ws.Cells.Item(22, 13).Value = -10
$ws.Cells.Item(22, 14).Value = -3404.5715
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Cells.Item(27, 8).Value = 2256.889
$ws.Cells.Item(27, 9).Value = 305
$ws.Cells.Item(27, 10).Value = 2814.5715
$ws.Cells.Item(27, 11).Value = 305
$ws.Cells.Item(27, 12).Value = 2814.5715
$ws.Cells.Item(27, 13).Value = -198
$ws.Cells.Item(27, 14).Value = -3028.5715
# Row 36: Campaign in the Membrane / Toadskin Jacket
$ws.Cells.Item(36, 8).Value = 25000
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = 25000
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 12).Value = 25000
$ws.Cells.Item(36, 14).Value = -26124
# Row 40: Best Served Toad / Toad Leather
$ws.Cells.Item(40, 8).Value = 27472.65
$ws.Cells.Item(40, 9).Value = 46701.184
$ws.Cells.Item(40, 10).Value = 3971.111
$ws.Cells.Item(40, 11).Value = 46701.184
$ws.Cells.Item(40, 12).Value = 3971.111
$ws.Cells.Item(40, 13).Value = -46565.184
$ws.Cells.Item(40, 14).Value = -4243.111
# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Cells.Item(55, 8).Value = 687.2727
$ws.Cells.Item(55, 9).Value = 592.5
$ws.Cells.Item(55, 10).Value = 940
$ws.Cells.Item(55, 11).Value = 592.5
$ws.Cells.Item(55, 12).Value = 940
$ws.Cells.Item(55, 13).Value = -419.5
$ws.Cells.Item(55, 14).Value = -1286
# Row 82: Trainin' the Neck / Dragon Leather
$ws.Cells.Item(82, 8).Value = 1430.6923
$ws.Cells.Item(82, 9).Value = 1133.1666
$ws.Cells.Item(82, 10).Value = 1685.7142
$ws.Cells.Item(82, 11).Value = 1133.1666
$ws.Cells.Item(82, 12).Value = 1685.7142
$ws.Cells.Item(82, 13).Value = -772.1666
$ws.Cells.Item(82, 14).Value = -2407.7142
# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Cells.Item(85, 8).Value = 1430.6923
$ws.Cells.Item(85, 9).Value = 1133.1666
$ws.Cells.Item(85, 10).Value = 1685.7142
$ws.Cells.Item(85, 11).Value = 1133.1666
$ws.Cells.Item(85, 12).Value = 1685.7142
$ws.Cells.Item(85, 13).Value = 114.8334
$ws.Cells.Item(85, 14).Value = -4181.7142
# Row 122: Hell on Leather / Gaja Leather
$ws.Cells.Item(122, 8).Value = 5294614.5
$ws.Cells.Item(122, 9).Value = 15874177
$ws.Cells.Item(122, 10).Value = 4832.7856
$ws.Cells.Item(122, 11).Value = 47622531
$ws.Cells.Item(122, 12).Value = 14498.3568
$ws.Cells.Item(122, 13).Value = -47620081
$ws.Cells.Item(122, 14).Value = -19398.3568
